$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 84 (Indice 83): F. Amager vs Esbjerg, copying the row-83 formatting ---
$ws.Range("A83:V83").Copy()
$ws.Range("A84").PasteSpecial(-4122)

$ws.Cells.Item(84, 1).Value = 83
$ws.Cells.Item(84, 2).Value = "denmark"
$ws.Cells.Item(84, 3).Value = "2nd-division"
$ws.Cells.Item(84, 4).Value = "2023-2024"
$ws.Cells.Item(84, 5).Value = 45235.54166666666
$ws.Cells.Item(84, 6).Value = "F. Amager"
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = "Esbjerg"
$ws.Cells.Item(84, 9).Value = 1
$ws.Cells.Item(84, 10).Value = 4.63
$ws.Cells.Item(84, 11).Value = "04/11/2023 01:13"
$ws.Cells.Item(84, 12).Value = 7.04
$ws.Cells.Item(84, 13).Value = "05/11/2023 12:59"
$ws.Cells.Item(84, 14).Value = 4.4
$ws.Cells.Item(84, 15).Value = "04/11/2023 01:13"
$ws.Cells.Item(84, 16).Value = 5.46
$ws.Cells.Item(84, 17).Value = "05/11/2023 12:59"
$ws.Cells.Item(84, 18).Value = 1.49
$ws.Cells.Item(84, 19).Value = "04/11/2023 01:13"
$ws.Cells.Item(84, 20).Value = 1.34
$ws.Cells.Item(84, 21).Value = "05/11/2023 12:59"
$ws.Cells.Item(84, 22).Value = "https://www.betexplorer.com/football/denmark/2nd-division/fremad-amager-esbjerg/Q5bFicm4/"

# --- Swap match data for rows that were reordered in the source scrape ---
# row 3 <= original row 4
$ws.Range("F3").Value = "Skive"
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = "Esbjerg"
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 3.83
$ws.Range("K3").Value = "03/08/2023 07:12"
$ws.Range("L3").Value = 4.58
$ws.Range("M3").Value = "04/08/2023 18:50"
$ws.Range("N3").Value = 3.73
$ws.Range("O3").Value = "03/08/2023 07:12"
$ws.Range("P3").Value = 4.05
$ws.Range("Q3").Value = "04/08/2023 18:50"
$ws.Range("R3").Value = 1.7
$ws.Range("S3").Value = "03/08/2023 07:12"
$ws.Range("T3").Value = 1.65
$ws.Range("U3").Value = "04/08/2023 18:50"
$ws.Range("V3").Value = "https://www.betexplorer.com/football/denmark/2nd-division/skive-esbjerg/tEGoB9XP/"

# row 4 <= original row 3
$ws.Range("F4").Value = "F. Amager"
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = "Hellerup"
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1.91
$ws.Range("K4").Value = "04/08/2023 17:43"
$ws.Range("L4").Value = 1.91
$ws.Range("M4").Value = "04/08/2023 17:43"
$ws.Range("N4").Value = 3.73
$ws.Range("O4").Value = "04/08/2023 17:43"
$ws.Range("P4").Value = 3.73
$ws.Range("Q4").Value = "04/08/2023 17:43"
$ws.Range("R4").Value = 3.56
$ws.Range("S4").Value = "04/08/2023 17:43"
$ws.Range("T4").Value = 3.56
$ws.Range("U4").Value = "04/08/2023 17:43"
$ws.Range("V4").Value = "https://www.betexplorer.com/football/denmark/2nd-division/fremad-amager-hellerup/2NHsCkIJ/"

# row 14 <= original row 15
$ws.Range("F14").Value = "Nykobing"
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = "Roskilde"
$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 2.04
$ws.Range("K14").Value = "17/08/2023 21:49"
$ws.Range("L14").Value = 2.38
$ws.Range("M14").Value = "18/08/2023 17:37"
$ws.Range("N14").Value = 3.5
$ws.Range("O14").Value = "17/08/2023 21:49"
$ws.Range("P14").Value = 3.81
$ws.Range("Q14").Value = "18/08/2023 17:34"
$ws.Range("R14").Value = 3.31
$ws.Range("S14").Value = "17/08/2023 21:49"
$ws.Range("T14").Value = 2.57
$ws.Range("U14").Value = "18/08/2023 17:37"
$ws.Range("V14").Value = "https://www.betexplorer.com/football/denmark/2nd-division/nykobing-roskilde/KtZ6G495/"

# row 15 <= original row 16
$ws.Range("F15").Value = "AB Copenhagen"
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = "Hellerup"
$ws.Range("I15").Value = 4
$ws.Range("J15").Value = 2.01
$ws.Range("K15").Value = "17/08/2023 07:12"
$ws.Range("L15").Value = 1.97
$ws.Range("M15").Value = "18/08/2023 18:50"
$ws.Range("N15").Value = 3.53
$ws.Range("O15").Value = "17/08/2023 07:12"
$ws.Range("P15").Value = 3.65
$ws.Range("Q15").Value = "18/08/2023 18:50"
$ws.Range("R15").Value = 3
$ws.Range("S15").Value = "17/08/2023 07:12"
$ws.Range("T15").Value = 3.46
$ws.Range("U15").Value = "18/08/2023 18:50"
$ws.Range("V15").Value = "https://www.betexplorer.com/football/denmark/2nd-division/ab-copenhagen-hellerup/86V2Hpfa/"

# row 16 <= original row 14
$ws.Range("F16").Value = "Skive"
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = "Aarhus Fremad"
$ws.Range("I16").Value = 2
$ws.Range("J16").Value = 4.46
$ws.Range("K16").Value = "17/08/2023 07:12"
$ws.Range("L16").Value = 4.82
$ws.Range("M16").Value = "18/08/2023 08:52"
$ws.Range("N16").Value = 3.89
$ws.Range("O16").Value = "17/08/2023 07:12"
$ws.Range("P16").Value = 4.12
$ws.Range("Q16").Value = "18/08/2023 17:02"
$ws.Range("R16").Value = 1.6
$ws.Range("S16").Value = "17/08/2023 07:12"
$ws.Range("T16").Value = 1.6
$ws.Range("U16").Value = "17/08/2023 07:12"
$ws.Range("V16").Value = "https://www.betexplorer.com/football/denmark/2nd-division/skive-aarhus-fremad/EqWbIQvg/"

# row 22 <= original row 23
$ws.Range("F22").Value = "Roskilde"
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = "Skive"
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 1.76
$ws.Range("K22").Value = "25/08/2023 02:12"
$ws.Range("L22").Value = 2.26
$ws.Range("M22").Value = "26/08/2023 13:47"
$ws.Range("N22").Value = 3.58
$ws.Range("O22").Value = "25/08/2023 02:12"
$ws.Range("P22").Value = 3.44
$ws.Range("Q22").Value = "26/08/2023 13:47"
$ws.Range("R22").Value = 3.75
$ws.Range("S22").Value = "25/08/2023 02:12"
$ws.Range("T22").Value = 2.96
$ws.Range("U22").Value = "26/08/2023 13:47"
$ws.Range("V22").Value = "https://www.betexplorer.com/football/denmark/2nd-division/roskilde-skive/SndYWPgn/"

# row 23 <= original row 22
$ws.Range("F23").Value = "Middelfart"
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = "F. Amager"
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 2.14
$ws.Range("K23").Value = "26/08/2023 10:12"
$ws.Range("L23").Value = 2
$ws.Range("M23").Value = "26/08/2023 13:25"
$ws.Range("N23").Value = 3.61
$ws.Range("O23").Value = "26/08/2023 10:12"
$ws.Range("P23").Value = 3.62
$ws.Range("Q23").Value = "26/08/2023 13:25"
$ws.Range("R23").Value = 2.91
$ws.Range("S23").Value = "26/08/2023 10:12"
$ws.Range("T23").Value = 3.39
$ws.Range("U23").Value = "26/08/2023 13:25"
$ws.Range("V23").Value = "https://www.betexplorer.com/football/denmark/2nd-division/middelfart-fremad-amager/M9exWq9h/"

# row 31 <= original row 33
$ws.Range("F31").Value = "FA 2000"
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = "Hellerup"
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 2.18
$ws.Range("K31").Value = "08/09/2023 08:12"
$ws.Range("L31").Value = 2.08
$ws.Range("M31").Value = "08/09/2023 17:30"
$ws.Range("N31").Value = 3.44
$ws.Range("O31").Value = "08/09/2023 08:12"
$ws.Range("P31").Value = 3.56
$ws.Range("Q31").Value = "08/09/2023 18:32"
$ws.Range("R31").Value = 2.95
$ws.Range("S31").Value = "08/09/2023 08:12"
$ws.Range("T31").Value = 3.02
$ws.Range("U31").Value = "08/09/2023 17:30"
$ws.Range("V31").Value = "https://www.betexplorer.com/football/denmark/2nd-division/frederiksberg-alliancen-2000-hellerup/0Gr4pqOo/"

# row 32 <= original row 31
$ws.Range("F32").Value = "Thisted FC"
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = "Skive"
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2.16
$ws.Range("K32").Value = "07/09/2023 07:12"
$ws.Range("L32").Value = 2.41
$ws.Range("M32").Value = "08/09/2023 18:19"
$ws.Range("N32").Value = 3.37
$ws.Range("O32").Value = "07/09/2023 07:12"
$ws.Range("P32").Value = 3.39
$ws.Range("Q32").Value = "08/09/2023 18:56"
$ws.Range("R32").Value = 2.9
$ws.Range("S32").Value = "07/09/2023 07:12"
$ws.Range("T32").Value = 2.77
$ws.Range("U32").Value = "08/09/2023 18:56"
$ws.Range("V32").Value = "https://www.betexplorer.com/football/denmark/2nd-division/thisted-fc-skive/MVWCrNhb/"

# row 33 <= original row 32
$ws.Range("F33").Value = "Roskilde"
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = "AB Copenhagen"
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 2.15
$ws.Range("K33").Value = "07/09/2023 07:12"
$ws.Range("L33").Value = 2.25
$ws.Range("M33").Value = "08/09/2023 18:58"
$ws.Range("N33").Value = 3.42
$ws.Range("O33").Value = "07/09/2023 07:12"
$ws.Range("P33").Value = 3.34
$ws.Range("Q33").Value = "08/09/2023 18:58"
$ws.Range("R33").Value = 2.79
$ws.Range("S33").Value = "07/09/2023 07:12"
$ws.Range("T33").Value = 3.05
$ws.Range("U33").Value = "08/09/2023 18:58"
$ws.Range("V33").Value = "https://www.betexplorer.com/football/denmark/2nd-division/roskilde-ab-copenhagen/v7s8q3wh/"

# row 51 <= original row 52
$ws.Range("F51").Value = "Brabrand"
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = "Skive"
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2.69
$ws.Range("K51").Value = "29/09/2023 01:13"
$ws.Range("L51").Value = 2.87
$ws.Range("M51").Value = "30/09/2023 12:03"
$ws.Range("N51").Value = 3.13
$ws.Range("O51").Value = "29/09/2023 01:13"
$ws.Range("P51").Value = 3.18
$ws.Range("Q51").Value = "30/09/2023 13:53"
$ws.Range("R51").Value = 2.42
$ws.Range("S51").Value = "29/09/2023 01:13"
$ws.Range("T51").Value = 2.45
$ws.Range("U51").Value = "30/09/2023 13:53"
$ws.Range("V51").Value = "https://www.betexplorer.com/football/denmark/2nd-division/brabrand-skive/IZaweirL/"

# row 52 <= original row 51
$ws.Range("F52").Value = "Roskilde"
$ws.Range("G52").Value = 3
$ws.Range("H52").Value = "Esbjerg"
$ws.Range("I52").Value = 3
$ws.Range("J52").Value = 3.56
$ws.Range("K52").Value = "29/09/2023 01:13"
$ws.Range("L52").Value = 3.81
$ws.Range("M52").Value = "30/09/2023 13:57"
$ws.Range("N52").Value = 3.63
$ws.Range("O52").Value = "29/09/2023 01:13"
$ws.Range("P52").Value = 3.76
$ws.Range("Q52").Value = "30/09/2023 12:01"
$ws.Range("R52").Value = 1.79
$ws.Range("S52").Value = "29/09/2023 01:13"
$ws.Range("T52").Value = 1.84
$ws.Range("U52").Value = "30/09/2023 13:57"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/denmark/2nd-division/roskilde-esbjerg/zBAYeXSE/"

# row 75 <= original row 76
$ws.Range("F75").Value = "Roskilde"
$ws.Range("G75").Value = 4
$ws.Range("H75").Value = "Brabrand"
$ws.Range("I75").Value = 4
$ws.Range("J75").Value = 1.56
$ws.Range("K75").Value = "27/10/2023 02:12"
$ws.Range("L75").Value = 1.48
$ws.Range("M75").Value = "27/10/2023 10:36"
$ws.Range("N75").Value = 3.92
$ws.Range("O75").Value = "27/10/2023 02:12"
$ws.Range("P75").Value = 4.27
$ws.Range("Q75").Value = "28/10/2023 12:03"
$ws.Range("R75").Value = 4.56
$ws.Range("S75").Value = "27/10/2023 02:12"
$ws.Range("T75").Value = 6.09
$ws.Range("U75").Value = "27/10/2023 10:36"
$ws.Range("V75").Value = "https://www.betexplorer.com/football/denmark/2nd-division/roskilde-brabrand/Wv7vjpCa/"

# row 76 <= original row 75
$ws.Range("F76").Value = "Aarhus Fremad"
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = "Middelfart"
$ws.Range("I76").Value = 2
$ws.Range("J76").Value = 1.73
$ws.Range("K76").Value = "27/10/2023 02:12"
$ws.Range("L76").Value = 1.82
$ws.Range("M76").Value = "28/10/2023 13:51"
$ws.Range("N76").Value = 3.76
$ws.Range("O76").Value = "27/10/2023 02:12"
$ws.Range("P76").Value = 3.73
$ws.Range("Q76").Value = "28/10/2023 13:51"
$ws.Range("R76").Value = 3.69
$ws.Range("S76").Value = "27/10/2023 02:12"
$ws.Range("T76").Value = 3.94
$ws.Range("U76").Value = "28/10/2023 13:51"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/denmark/2nd-division/aarhus-fremad-middelfart/Iy3ziQdg/"

# row 77 <= original row 78
$ws.Range("F77").Value = "Thisted FC"
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = "AB Copenhagen"
$ws.Range("I77").Value = 1
$ws.Range("J77").Value = 2.52
$ws.Range("K77").Value = "27/10/2023 03:12"
$ws.Range("L77").Value = 2.32
$ws.Range("M77").Value = "28/10/2023 14:59"
$ws.Range("N77").Value = 3.33
$ws.Range("O77").Value = "27/10/2023 03:12"
$ws.Range("P77").Value = 3.54
$ws.Range("Q77").Value = "28/10/2023 14:59"
$ws.Range("R77").Value = 2.4
$ws.Range("S77").Value = "27/10/2023 03:12"
$ws.Range("T77").Value = 2.79
$ws.Range("U77").Value = "28/10/2023 14:59"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/denmark/2nd-division/thisted-fc-ab-copenhagen/hfemlOtC/"

# row 78 <= original row 77
$ws.Range("F78").Value = "Nykobing"
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = "Esbjerg"
$ws.Range("I78").Value = 3
$ws.Range("J78").Value = 3.87
$ws.Range("K78").Value = "27/10/2023 03:12"
$ws.Range("L78").Value = 4.51
$ws.Range("M78").Value = "28/10/2023 14:34"
$ws.Range("N78").Value = 4.01
$ws.Range("O78").Value = "27/10/2023 03:12"
$ws.Range("P78").Value = 4.37
$ws.Range("Q78").Value = "28/10/2023 14:34"
$ws.Range("R78").Value = 1.65
$ws.Range("S78").Value = "27/10/2023 03:12"
$ws.Range("T78").Value = 1.61
$ws.Range("U78").Value = "28/10/2023 09:39"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/denmark/2nd-division/nykobing-esbjerg/E3dqk4R5/"

# row 80 <= original row 82
$ws.Range("F80").Value = "Middelfart"
$ws.Range("G80").Value = 2
$ws.Range("H80").Value = "FA 2000"
$ws.Range("I80").Value = 1
$ws.Range("J80").Value = 1.5
$ws.Range("K80").Value = "03/11/2023 02:12"
$ws.Range("L80").Value = 1.58
$ws.Range("M80").Value = "04/11/2023 13:52"
$ws.Range("N80").Value = 4.13
$ws.Range("O80").Value = "03/11/2023 02:12"
$ws.Range("P80").Value = 4.06
$ws.Range("Q80").Value = "04/11/2023 13:52"
$ws.Range("R80").Value = 5.07
$ws.Range("S80").Value = "03/11/2023 02:12"
$ws.Range("T80").Value = 5.25
$ws.Range("U80").Value = "04/11/2023 13:52"
$ws.Range("V80").Value = "https://www.betexplorer.com/football/denmark/2nd-division/middelfart-frederiksberg-alliancen-2000/nFaBhwYc/"

# row 82 <= original row 80
$ws.Range("F82").Value = "Brabrand"
$ws.Range("G82").Value = 2
$ws.Range("H82").Value = "Nykobing"
$ws.Range("I82").Value = 1
$ws.Range("J82").Value = 3.27
$ws.Range("K82").Value = "03/11/2023 02:12"
$ws.Range("L82").Value = 3.29
$ws.Range("M82").Value = "04/11/2023 13:47"
$ws.Range("N82").Value = 3.38
$ws.Range("O82").Value = "03/11/2023 02:12"
$ws.Range("P82").Value = 3.54
$ws.Range("Q82").Value = "04/11/2023 13:47"
$ws.Range("R82").Value = 1.98
$ws.Range("S82").Value = "03/11/2023 02:12"
$ws.Range("T82").Value = 2.06
$ws.Range("U82").Value = "04/11/2023 13:47"
$ws.Range("V82").Value = "https://www.betexplorer.com/football/denmark/2nd-division/brabrand-nykobing/vkd3fa3o/"

